$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.077.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "'3.461.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'576.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'159.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'3.460.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'4.055.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "'0.135"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'28.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "'65.052.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "'3.485.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'14.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'380.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "'0.558"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").Value = "'72.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "'0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'10.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.42%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +4.90%  "
$ws.Range("D31").Value = "'6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "'7.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.84%  "
$ws.Range("E35").Value = "  +10.25%  "
$ws.Range("D36").Value = "'161.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'1.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").Value = "'0.0786"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'27.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.25%  "
$ws.Range("D41").Value = "'2.916.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'6.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.27%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0320"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'43.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "'0.781"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").Value = "'25.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.25%  "
$ws.Range("D47").Value = "'321.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.53%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'0.880"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
